# Updated cryptos list on Mon Apr  8 19:18:48 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed by Excel as a plain
# number (single decimal point, no letters) - force them to stay text so
# they round-trip exactly like the original inline-string cells.
$numericRiskCells = @("D5","D6","D8","D12","D15","D21","D22","D23","D25","D28","D29","D32","D33","D34","D36","D41","D44","D45","D51")
foreach ($addr in $numericRiskCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "71.918.56"
$ws.Range("E2").Value = "  +3.77%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.694.11"
$ws.Range("E3").Value = "  +9.24%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.06%  "

# Row 5 - BNB
$ws.Range("D5").Value = "589.32"
$ws.Range("E5").Value = "  +1.64%  "

# Row 6 - Solana
$ws.Range("D6").Value = "180.32"
$ws.Range("E6").Value = "  +1.16%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.687.63"

# Row 8 - XRP
$ws.Range("D8").Value = "0.625"
$ws.Range("E8").Value = "  +5.40%  "

# Row 9 - USDC
$ws.Range("E9").Value = "  +0.07%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +2.95%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  +4.77%  "

# Row 12 - Avalanche
$ws.Range("D12").Value = "50.09"
$ws.Range("E12").Value = "  +3.71%  "

# Row 13 - ShibaInu
$ws.Range("E13").Value = "  +1.59%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "4.287.96"
$ws.Range("E14").Value = "  +9.34%  "

# Row 15 - BitcoinCash
$ws.Range("D15").Value = "685.24"
$ws.Range("E15").Value = "  +0.20%  "

# Row 16 - Polkadot
$ws.Range("E16").Value = "  +4.62%  "

# Row 17 / Row 18 - WrappedEther and WrappedBTC swapped ranking order
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "72.043.32"
$ws.Range("E17").Value = "  +3.79%  "

$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.690.63"
$ws.Range("E18").Value = "  +9.23%  "

# Row 19 - TRON
$ws.Range("E19").Value = "  +2.05%  "

# Row 20 - Chainlink
$ws.Range("E20").Value = "  +3.62%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "11.67"
$ws.Range("E21").Value = "  +3.62%  "

# Row 22 - Polygon
$ws.Range("D22").Value = "0.943"
$ws.Range("E22").Value = "  +3.94%  "

# Row 23 - Toncoin
$ws.Range("D23").Value = "6.08"
$ws.Range("E23").Value = "  +13.79%  "

# Row 24 - InternetComputer(DFINITY)
$ws.Range("E24").Value = "  +3.94%  "

# Row 25 - Litecoin
$ws.Range("D25").Value = "103.76"
$ws.Range("E25").Value = "  +2.61%  "

# Row 26 - PancakeSwap
$ws.Range("E26").Value = "  +4.42%  "

# Row 27 - ImmutableX
$ws.Range("E27").Value = "  +6.32%  "

# Row 28 - RenderToken
$ws.Range("D28").Value = "10.24"
$ws.Range("E28").Value = "  +5.58%  "

# Row 29 - EthereumClassic
$ws.Range("D29").Value = "35.61"
$ws.Range("E29").Value = "  +6.58%  "

# Row 30 - Filecoin
$ws.Range("E30").Value = "  +5.74%  "

# Row 31 - NEARProtocol
$ws.Range("E31").Value = "  +7.23%  "

# Row 32 - dogwifhat
$ws.Range("D32").Value = "4.21"
$ws.Range("E32").Value = "  +10.10%  "

# Row 33 - Bittensor
$ws.Range("D33").Value = "579.11"
$ws.Range("E33").Value = "  +5.51%  "

# Row 34 - Cosmos
$ws.Range("D34").Value = "11.34"
$ws.Range("E34").Value = "  +3.00%  "

# Row 35 - Hedera
$ws.Range("E35").Value = "  +4.22%  "

# Row 36 - OKB
$ws.Range("D36").Value = "59.91"
$ws.Range("E36").Value = "  +3.67%  "

# Row 37 - Maker
$ws.Range("D37").Value = "3.764.02"
$ws.Range("E37").Value = "  +4.51%  "

# Row 38 - Dai
$ws.Range("E38").Value = "  +0.03%  "

# Row 39 - Kaspa
$ws.Range("E39").Value = "  +3.29%  "

# Row 40 - PEPE (subscript-3 character in the price)
$ws.Range("D40").Value = "0.0₃0775"
$ws.Range("E40").Value = "  +4.95%  "

# Row 41 - InjectiveProtocol
$ws.Range("D41").Value = "35.61"
$ws.Range("E41").Value = "  +1.00%  "

# Row 42 - Stacks
$ws.Range("E42").Value = "  +5.18%  "

# Row 43 - Fetch.AI
$ws.Range("E43").Value = "  +3.66%  "

# Row 44 - VeChain
$ws.Range("D44").Value = "0.0463"
$ws.Range("E44").Value = "  +9.17%  "

# Row 45 - TheGraph
$ws.Range("D45").Value = "0.349"
$ws.Range("E45").Value = "  +4.29%  "

# Row 46 - ApeXProtocol
$ws.Range("E46").Value = "  +0.09%  "

# Row 47 - ThetaToken
$ws.Range("E47").Value = "  +6.89%  "

# Row 48 - Stellar (unchanged)

# Row 49 - Mantle
$ws.Range("E49").Value = "  +3.82%  "

# Row 50 - FirstDigitalUSD
$ws.Range("E50").Value = "  -0.16%  "

# Row 51 - Monero
$ws.Range("D51").Value = "133.96"
$ws.Range("E51").Value = "  +3.45%  "
